# Daily BP terminal gate pricing (TGP) refresh.
#
# The workbook carries a rolling two-day window per terminal: the newest
# "Effective Date" row and the prior day's row. This update rolls the
# window forward by one day (old "today" rows become "yesterday" rows
# with their same figures, and brand-new figures are published for the
# new "today"), across every state section on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => { Column => new value }. Column "A" is the Excel date serial
# (the cell is already date-formatted, so assigning the serial number
# renders correctly); D/E/F/G are Diesel/ULP/PULP/e10 prices.
$updates = [ordered]@{
    8  = @{ "A"=46038; "D"=154.25; "E"=150.08; "F"=160.08; "G"=150.1 }
    9  = @{ "A"=46038; "D"=154.25; "E"=150.08; "F"=160.08; "G"=150.1 }
    10 = @{ "A"=46038; "D"=155.62; "E"=152.53; "F"=162.53; "G"=152.93 }
    11 = @{ "A"=46037; "D"=152.67; "E"=149.22; "F"=159.22; "G"=149.23 }
    12 = @{ "A"=46037; "D"=152.67; "E"=149.22; "F"=159.22; "G"=149.23 }
    13 = @{ "A"=46037; "D"=154.32; "E"=151.72; "F"=161.72; "G"=152.12 }

    17 = @{ "A"=46038; "D"=159.38; "E"=155.16; "F"=165.16 }
    18 = @{ "A"=46037; "D"=158.11; "E"=154.37; "F"=164.37 }

    22 = @{ "A"=46038; "D"=155.66; "E"=152.09; "F"=161.69; "G"=153.17 }
    23 = @{ "A"=46038; "D"=160.6;  "E"=157.87; "F"=167.87 }
    24 = @{ "A"=46038; "D"=160.75; "E"=158.53; "F"=168.53 }
    25 = @{ "A"=46038; "D"=160.73; "E"=158.05; "F"=168.05; "G"=158.18 }
    26 = @{ "A"=46038; "D"=160.35; "E"=159.67; "F"=169.67 }
    27 = @{ "A"=46037; "D"=154.08; "E"=151.23; "F"=160.83; "G"=152.3 }
    28 = @{ "A"=46037; "D"=159.29; "E"=157.06; "F"=167.06 }
    29 = @{ "A"=46037; "D"=159.45; "E"=157.7;  "F"=167.7 }
    30 = @{ "A"=46037; "D"=159.43; "E"=157.23; "F"=167.23; "G"=157.36 }
    31 = @{ "A"=46037; "D"=159.06; "E"=158.84; "F"=168.84 }

    35 = @{ "A"=46038; "D"=154.06; "E"=149.51; "F"=158.51 }
    36 = @{ "A"=46037; "D"=152.21; "E"=148.7;  "F"=157.7 }

    40 = @{ "A"=46038; "D"=160.16; "E"=157.82; "F"=167.82 }
    41 = @{ "A"=46038; "D"=159.88; "E"=158.24; "F"=168.24 }
    42 = @{ "A"=46037; "D"=158.84; "E"=157.1;  "F"=167.1 }
    43 = @{ "A"=46037; "D"=158.56; "E"=157.52; "F"=167.52 }

    47 = @{ "A"=46038; "D"=153.09; "E"=150.47; "F"=160.47 }
    48 = @{ "A"=46038; "D"=152.72; "E"=150.41; "F"=160.41 }
    49 = @{ "A"=46037; "D"=152.34; "E"=150.18; "F"=160.18 }
    50 = @{ "A"=46037; "D"=151.97; "E"=150.11; "F"=160.11 }

    54 = @{ "A"=46038; "D"=169.33; "E"=165.24; "F"=175.24 }
    55 = @{ "A"=46038; "D"=161.84; "E"=163.55; "F"=173.55 }
    56 = @{ "A"=46038; "D"=158.54 }
    57 = @{ "A"=46038; "D"=159.1;  "E"=157.97 }
    58 = @{ "A"=46038; "D"=154.87; "E"=153.87; "F"=163.87 }
    59 = @{ "A"=46038; "D"=161.72; "E"=163.52 }
    60 = @{ "A"=46037; "D"=168.03; "E"=164.32; "F"=174.32 }
    61 = @{ "A"=46037; "D"=160.49; "E"=162.44; "F"=172.44 }
    62 = @{ "A"=46037; "D"=156.91 }
    63 = @{ "A"=46037; "D"=157.51; "E"=156.86 }
    64 = @{ "A"=46037; "D"=153.28; "E"=152.76; "F"=162.76 }
    65 = @{ "A"=46037; "D"=160.46; "E"=162.65 }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value = $cells[$col]
    }
}
